$wb = $excel.ActiveWorkbook

$newText = "Here is a list of everyone that has been added already"

# "Create Full Input" sheet: update the two "Monday Night Football Squad"
# cells with the new text, then leave the selection on E2 (matches the
# final state captured in the workbook).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E2").Value = $newText
$ws1.Range("E3").Value = $newText
$ws1.Columns.Item(5).ColumnWidth = 48.7
$ws1.Range("E2").Select() | Out-Null

# "List Remove" sheet: selection moved to I26, and it is no longer the
# active sheet.
$ws13 = $wb.Worksheets.Item(13)
$ws13.Range("I26").Select() | Out-Null

# "Login Full" sheet: update the matching two cells, then make it the
# active sheet with the selection on E8.
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("C2").Value = $newText
$ws7.Range("C3").Value = $newText
$ws7.Columns.Item(3).ColumnWidth = 48.7
$ws7.Activate() | Out-Null
$ws7.Range("E8").Select() | Out-Null
